# Apply "last report 19-02-25" updates to the Route Cost RSO workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Route")

# RSO 02 (row 7) route cost updated 130 -> 170.
$ws.Range("D7").Value = 170

# RSO 03/04 section: row 9 (Liton Ray) 130 -> 200.
$ws.Range("D9").Value = 200

# Row 10: SR changed from "Ripon Mondal" to "Arman Hossen", cost 200 -> 150.
$ws.Range("C10").Value = "Arman Hossen"
$ws.Range("D10").Value = 150

# Header date cell L3: was a date serial (45690), now a free-text label "19/2025".
$ws.Range("L3").Value = "19/2025"
